$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-51: row, Coin (B), Link (C), Price (D), Volume(1h) (E).
# A new coin "OKB" is inserted at row 9, pushing Dogecoin..Quant down by one
# row each; "WOONetwork" (previously the last row, 51) drops off the bottom
# of the list. The D/E columns are refreshed with the latest market data for
# every row (2-51).
$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.520.31', '  +0.98%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.971.87', '  +3.66%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.006', '  +0.24%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '327.23', '  +0.36%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.006', '  +0.38%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4663', '  +0.35%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3913', '  -0.11%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '46.18', '  +0.11%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07934', '  +0.60%  '),
    @(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.9878', '  -0.23%  '),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '22.75', '  +4.41%  '),
    @(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.020.79', '  +6.04%  '),
    @(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.178', '  +1.45%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.801', '  +1.02%  '),
    @(16, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07126', '  +1.94%  '),
    @(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '87.73', '  -0.58%  '),
    @(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.008', '  +0.40%  '),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009914', '  -0.68%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.22', '  +0.69%  '),
    @(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.007', '  +0.55%  '),
    @(22, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.591.84', '  +1.21%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.542', '  +4.54%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.13', '  +0.48%  '),
    @(25, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.248.82', '  +5.77%  '),
    @(26, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.113', '  +0.02%  '),
    @(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '158.82', '  +1.62%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.53', '  +0.60%  '),
    @(29, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.791', '  -2.99%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '119.71', '  +0.89%  '),
    @(31, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.888', '  +0.22%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09424', '  +1.04%  '),
    @(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.8783', '  -2.66%  '),
    @(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.232', '  -0.49%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.316', '  -0.71%  '),
    @(36, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.151', '  -1.27%  '),
    @(37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05806', '  +0.57%  '),
    @(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.165', '  -1.54%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02102', '  +0.69%  '),
    @(40, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '7.701', '  -0.26%  '),
    @(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5706', '  +0.05%  '),
    @(42, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1795', '  +0.24%  '),
    @(43, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '9.623', '  -0.94%  '),
    @(44, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.755', '  +6.96%  '),
    @(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '11.78', '  -1.82%  '),
    @(46, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5327', '  -0.56%  '),
    @(47, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.000002735', '  +45.87%  '),
    @(48, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.140', '  -1.45%  '),
    @(49, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06926', '  -1.33%  '),
    @(50, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.825', '  -1.54%  '),
    @(51, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '113.40', '  +0.09%  ')

)

function Set-TextCell($cell, $value) {
    # D-column values are plain digit/dot strings (e.g. "327.23", "113.40")
    # that Excel's type-inference would otherwise silently convert to
    # numbers (losing the exact textual representation, e.g. trailing
    # zeros). Force the cell to Text, write the value, then drop back to
    # the default "Normal" style so no stray number-format sticks around.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    Set-TextCell $ws.Cells.Item($r, 4) $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
